$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.248583466539117
$ws.Range("C2").Value = 0.5135445598622255
$ws.Range("D2").Value = 0.02697253705651548
$ws.Range("E2").Value = 0.4183939171993387
$ws.Range("F2").Value = 1.40831938304018
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.8201472520215134

$ws.Range("B3").Value = 1.109004199223932
$ws.Range("C3").Value = 0.4494248925602164
$ws.Range("D3").Value = 0.02650423200020668
$ws.Range("E3").Value = 0.3645020341276677
$ws.Range("F3").Value = 1.3348920828372
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.793131367362335

$ws.Range("B4").Value = 1.023753385397583
$ws.Range("C4").Value = 0.4101930759799188
$ws.Range("D4").Value = 0.0262425623976732
$ws.Range("E4").Value = 0.3315723959174477
$ws.Range("F4").Value = 1.290923878185708
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.7771895207634003

$ws.Range("B5").Value = 0.9891235662772715
$ws.Range("C5").Value = 0.3942382203240413
$ws.Range("D5").Value = 0.02614226112430273
$ws.Range("E5").Value = 0.3181899391099563
$ws.Range("F5").Value = 1.273281674887869
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.7708526310546588

$ws.Range("B6").Value = 0.9833799061866557
$ws.Range("C6").Value = 0.3915908223520432
$ws.Range("D6").Value = 0.02612598349378459
$ws.Range("E6").Value = 0.3159698997834681
$ws.Range("F6").Value = 1.270368660628677
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.7698099468155704

$ws.Range("B7").Value = 1.023285911619382
$ws.Range("C7").Value = 0.4099777748626821
$ws.Range("D7").Value = 0.02624118429420008
$ws.Range("E7").Value = 0.3313917718563317
$ws.Range("F7").Value = 1.290684842283625
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.7771034170905438

$ws.Range("B8").Value = 1.200361005336902
$ws.Range("C8").Value = 0.4914060859751999
$ws.Range("D8").Value = 0.0268055902312696
$ws.Range("E8").Value = 0.3997765224959551
$ws.Range("F8").Value = 1.382766553573703
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.8106963311955866

$ws.Range("B9").Value = 1.551340953639112
$ws.Range("C9").Value = 0.65229212190917
$ws.Range("D9").Value = 0.02812564034314846
$ws.Range("E9").Value = 0.5353223072584257
$ws.Range("F9").Value = 1.572458603900571
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.8818331445848173

$ws.Range("B10").Value = 1.811729519971095
$ws.Range("C10").Value = 0.7713974297906248
$ws.Range("D10").Value = 0.02923668437058069
$ws.Range("E10").Value = 0.6360383629342294
$ws.Range("F10").Value = 1.717764074275038
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.9374986201205786

$ws.Range("B11").Value = 1.930791820376669
$ws.Range("C11").Value = 0.8258153070449907
$ws.Range("D11").Value = 0.02977525142903659
$ws.Range("E11").Value = 0.6821581709987896
$ws.Range("F11").Value = 1.785240746264435
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.9636035652325461

$ws.Range("B12").Value = 1.975969605603495
$ws.Range("C12").Value = 0.8464588475313803
$ws.Range("D12").Value = 0.02998416302251883
$ws.Range("E12").Value = 0.6996706982716745
$ws.Range("F12").Value = 1.810997008541079
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.9736047686078706

$ws.Range("B13").Value = 1.966235631998245
$ws.Range("C13").Value = 0.8420112170801985
$ws.Range("D13").Value = 0.02993894597989311
$ws.Range("E13").Value = 0.6958968596913166
$ws.Range("F13").Value = 1.805440747184122
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.9714456257839004

$ws.Range("B14").Value = 1.934506768806273
$ws.Range("C14").Value = 0.8275129128982712
$ws.Range("D14").Value = 0.02979233798153302
$ws.Range("E14").Value = 0.6835979510093608
$ws.Range("F14").Value = 1.787355597873272
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.9644240297522373

$ws.Range("B15").Value = 1.915083949766824
$ws.Range("C15").Value = 0.8186371398834922
$ws.Range("D15").Value = 0.02970318935458494
$ws.Range("E15").Value = 0.6760708892048086
$ws.Range("F15").Value = 1.776304720054412
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.9601382827241025

$ws.Range("B16").Value = 1.803961113385128
$ws.Range("C16").Value = 0.7678460767689899
$ws.Range("D16").Value = 0.02920217189517871
$ws.Range("E16").Value = 0.633030799230653
$ws.Range("F16").Value = 1.713382555926671
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.9358086273887807

$ws.Range("B17").Value = 1.735949755135266
$ws.Range("C17").Value = 0.7367498156958163
$ws.Range("D17").Value = 0.02890344601022576
$ws.Range("E17").Value = 0.6067078802236949
$ws.Range("F17").Value = 1.675138899022784
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.921085915907824

$ws.Range("B18").Value = 1.696888627717556
$ws.Range("C18").Value = 0.7188861426147923
$ws.Range("D18").Value = 0.02873473454552311
$ws.Range("E18").Value = 0.591595943341261
$ws.Range("F18").Value = 1.653271332854302
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.9126912477417193

$ws.Range("B19").Value = 1.683672915268005
$ws.Range("C19").Value = 0.7128415201045186
$ws.Range("D19").Value = 0.02867813942287967
$ws.Range("E19").Value = 0.5864840339517627
$ws.Range("F19").Value = 1.645889343899967
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.9098614716022126

$ws.Range("B20").Value = 1.743183733073693
$ws.Range("C20").Value = 0.7400577619174555
$ws.Range("D20").Value = 0.02893492280009013
$ws.Range("E20").Value = 0.6095070390605883
$ws.Range("F20").Value = 1.679196586759787
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.9226455484459137

$ws.Range("B21").Value = 1.943823792032958
$ws.Range("C21").Value = 0.8317703991930898
$ws.Range("D21").Value = 0.02983526385366986
$ws.Range("E21").Value = 0.6872091022820541
$ws.Range("F21").Value = 1.792662048340901
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.9664832725464976

$ws.Range("B22").Value = 2.075488690337409
$ws.Range("C22").Value = 0.8919249609699023
$ws.Range("D22").Value = 0.03045274161284794
$ws.Range("E22").Value = 0.738273492567771
$ws.Range("F22").Value = 1.868012203187078
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.9958104620327504

$ws.Range("B23").Value = 2.005166434535226
$ws.Range("C23").Value = 0.859798752631491
$ws.Range("D23").Value = 0.03012045635393434
$ws.Range("E23").Value = 0.7109922513340479
$ws.Range("F23").Value = 1.827685029803433
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.9800949631544853

$ws.Range("B24").Value = 1.739913128867101
$ws.Range("C24").Value = 0.7385621960838762
$ws.Range("D24").Value = 0.02892068272055326
$ws.Range("E24").Value = 0.6082414728318923
$ws.Range("F24").Value = 1.677361735266004
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.9219402219909512

$ws.Range("B25").Value = 1.455966638430652
$ws.Range("C25").Value = 0.6086220630136268
$ws.Range("D25").Value = 0.02774455462341763
$ws.Range("E25").Value = 0.4984734630056238
$ws.Range("F25").Value = 1.520126681014943
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.8620056120074793
